# Fill in the first week (rows 5-13, columns H:M = Sat..Fri) of the
# timesheet with the employee's reported hours. Row 14 (TOTAL), the
# N column sub-totals, and the P16:P25 "Pay Period total" column are all
# driven by formulas already present in the sheet, so they recalculate
# automatically once the inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns H..M correspond to Sat, Mon, Tue, Wed, Thu, Fri (row 4 headers).
$columns = @("H", "I", "J", "K", "L", "M")

# Ordered list of (row, Sat, Mon, Tue, Wed, Thu, Fri) hours for the first
# week's nine leave/work categories (rows 5-13).
$weekRows = @(
    , (5,  0, 2, 3, 4, 5, 6)
    , (6,  2, 2, 3, 4, 5, 6)
    , (7,  3, 2, 3, 4, 5, 6)
    , (8,  4, 0, 0, 0, 0, 0)
    , (9,  5, 0, 0, 0, 0, 0)
    , (10, 6, 0, 0, 0, 0, 0)
    , (11, 7, 0, 0, 0, 0, 0)
    , (12, 8, 0, 0, 0, 0, 0)
    , (13, 9, 0, 0, 0, 0, 0)
)

foreach ($entry in $weekRows) {
    $row = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $entry[$i + 1]
    }
}

# Mirror the author's recorded cursor/selection move to cell K11.
$ws.Range("K11").Select()
